# Apply scheduled-runner profit/price updates to the Siren_Profits workbook.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets per-row
# currentAveragePrice / LevePrice / LeveProfit corrections.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1146.55
$ws.Range("I19").Value = 676.5
$ws.Range("K19").Value = 676.5
$ws.Range("M19").Value = -501.5
# Row 43
$ws.Range("H43").Value = 3888.0557
$ws.Range("I43").Value = 2750
$ws.Range("K43").Value = 2750
$ws.Range("M43").Value = -2681
# Row 88
$ws.Range("H88").Value = 611.25
$ws.Range("I88").Value = 999
$ws.Range("K88").Value = 999
$ws.Range("M88").Value = -593
# Row 91
$ws.Range("H91").Value = 611.25
$ws.Range("I91").Value = 999
$ws.Range("K91").Value = 999
$ws.Range("M91").Value = 405
# Row 105
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
# Row 112
$ws.Range("H112").Value = 101367.09
$ws.Range("J112").Value = 368332.66
$ws.Range("L112").Value = 1104997.98
$ws.Range("N112").Value = -1107213.98
# Row 132
$ws.Range("H132").Value = 3521.738
$ws.Range("I132").Value = 3774.4412
$ws.Range("K132").Value = 11323.3236
$ws.Range("M132").Value = -8793.3236
# Row 138
$ws.Range("H138").Value = 4607.0586
$ws.Range("J138").Value = 5156.384
$ws.Range("L138").Value = 15469.152
$ws.Range("N138").Value = -25749.152

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2738.068
$ws.Range("I32").Value = 2670.195
$ws.Range("J32").Value = 3665.6667
$ws.Range("K32").Value = 2670.195
$ws.Range("L32").Value = 3665.6667
$ws.Range("M32").Value = -2383.195
$ws.Range("N32").Value = -4239.6667
# Row 102
$ws.Range("H102").Value = 7306.609
$ws.Range("I102").Value = 3110.7144
$ws.Range("K102").Value = 3110.7144
$ws.Range("M102").Value = -1488.7144
# Row 110
$ws.Range("H110").Value = 1500
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 545

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 70956.39999999999
$ws.Range("I82").Value = 27499.5
$ws.Range("J82").Value = 99927.664
$ws.Range("K82").Value = 27499.5
$ws.Range("L82").Value = 99927.664
$ws.Range("M82").Value = -27116.5
$ws.Range("N82").Value = -100693.664
# Row 85
$ws.Range("H85").Value = 70956.39999999999
$ws.Range("I85").Value = 27499.5
$ws.Range("J85").Value = 99927.664
$ws.Range("K85").Value = 27499.5
$ws.Range("L85").Value = 99927.664
$ws.Range("M85").Value = -26173.5
$ws.Range("N85").Value = -102579.664
# Row 97
$ws.Range("H97").Value = 24998.666
$ws.Range("I97").Value = 14997
$ws.Range("J97").Value = 29999.5
$ws.Range("K97").Value = 14997
$ws.Range("L97").Value = 29999.5
$ws.Range("M97").Value = -14006
$ws.Range("N97").Value = -31981.5
# Row 107
$ws.Range("H107").Value = 8533
$ws.Range("I107").Value = 8533
$ws.Range("K107").Value = 8533
$ws.Range("M107").Value = -6613
# Row 134
$ws.Range("H134").Value = 6085.75
$ws.Range("I134").Value = 6329.593
$ws.Range("J134").Value = 4769
$ws.Range("K134").Value = 18988.779
$ws.Range("L134").Value = 14307
$ws.Range("M134").Value = -16453.779
$ws.Range("N134").Value = -19377

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7178.3335
$ws.Range("I31").Value = 7239.143
$ws.Range("J31").Value = 7125.125
$ws.Range("K31").Value = 7239.143
$ws.Range("L31").Value = 7125.125
$ws.Range("M31").Value = -6944.143
$ws.Range("N31").Value = -7715.125
# Row 34
$ws.Range("H34").Value = 7178.3335
$ws.Range("I34").Value = 7239.143
$ws.Range("J34").Value = 7125.125
$ws.Range("K34").Value = 7239.143
$ws.Range("L34").Value = 7125.125
$ws.Range("M34").Value = -7037.143
$ws.Range("N34").Value = -7529.125
# Row 132
$ws.Range("H132").Value = 23438.48
$ws.Range("I132").Value = 8218.15
$ws.Range("K132").Value = 24654.45
$ws.Range("M132").Value = -22124.45

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 36608.94
$ws.Range("J97").Value = 1662.7
$ws.Range("L97").Value = 4988.1
$ws.Range("N97").Value = -5980.1
# Row 131
$ws.Range("H131").Value = 66668470
$ws.Range("I131").Value = 500000960
$ws.Range("J131").Value = 1932.8462
$ws.Range("K131").Value = 1500002880
$ws.Range("L131").Value = 5798.5386
$ws.Range("M131").Value = -1499997840
$ws.Range("N131").Value = -15878.5386

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 35388.5
$ws.Range("J52").Value = 35388.5
$ws.Range("L52").Value = 35388.5
$ws.Range("N52").Value = -35906.5
# Row 97
$ws.Range("H97").Value = 10633.913
$ws.Range("I97").Value = 11749.421
$ws.Range("J97").Value = 5335.25
$ws.Range("K97").Value = 11749.421
$ws.Range("L97").Value = 5335.25
$ws.Range("M97").Value = -11253.421
$ws.Range("N97").Value = -6327.25
# Row 122
$ws.Range("H122").Value = 11960.471
$ws.Range("I122").Value = 10473.556
$ws.Range("K122").Value = 31420.668
$ws.Range("M122").Value = -28970.668
# Row 132
$ws.Range("H132").Value = 3188.5908
$ws.Range("I132").Value = 2016.5834
$ws.Range("J132").Value = 4595
$ws.Range("K132").Value = 6049.7502
$ws.Range("L132").Value = 13785
$ws.Range("M132").Value = -3519.7502
$ws.Range("N132").Value = -18845

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 23094.6
$ws.Range("I81").Value = 29370.857
$ws.Range("J81").Value = 8450
$ws.Range("K81").Value = 58741.714
$ws.Range("L81").Value = 16900
$ws.Range("M81").Value = -57680.714
$ws.Range("N81").Value = -19022
# Row 84
$ws.Range("H84").Value = 23094.6
$ws.Range("I84").Value = 29370.857
$ws.Range("J84").Value = 8450
$ws.Range("K84").Value = 293708.57
$ws.Range("L84").Value = 84500
$ws.Range("M84").Value = -288404.57
$ws.Range("N84").Value = -95108
# Row 126
$ws.Range("H126").Value = 19338.44
$ws.Range("I126").Value = 23156.105
$ws.Range("J126").Value = 7249.1665
$ws.Range("K126").Value = 69468.315
$ws.Range("L126").Value = 21747.4995
$ws.Range("M126").Value = -66998.315
$ws.Range("N126").Value = -26687.4995
# Row 132
$ws.Range("H132").Value = 24870.117
$ws.Range("I132").Value = 26186.8
$ws.Range("K132").Value = 78560.39999999999
$ws.Range("M132").Value = -76030.39999999999
